$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105; this shifts the existing rows 105-129 down to 106-130
$ws.Rows.Item(105).Insert()

$newRow = 105

$ws.Cells.Item($newRow, 1).Value = 1
$ws.Cells.Item($newRow, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($newRow, 3).Value = "Arica y Parinacota"

$dCell = $ws.Cells.Item($newRow, 4)
$dCell.Value = 45180
$dCell.NumberFormat = $ws.Cells.Item(106, 4).NumberFormat

$ws.Cells.Item($newRow, 5).Value = 15
$ws.Cells.Item($newRow, 6).Value = 100112040
$ws.Cells.Item($newRow, 7).Value = "Cilantro"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 270
$ws.Cells.Item($newRow, 11).Value = 900
$ws.Cells.Item($newRow, 12).Value = 1000
$ws.Cells.Item($newRow, 13).Value = 950
$ws.Cells.Item($newRow, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item($newRow, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($newRow, 16).Value = 475
$ws.Cells.Item($newRow, 17).Value = 2
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
